$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-03 Monday" "2024-06-04 Tuesday"

Replace-Text "20÷5=4, 0" "83÷8=10, 3"
Replace-Text "11÷2=5, 1" "35÷5=7, 0"
Replace-Text "96÷3=32, 0" "28÷2=14, 0"
Replace-Text "15÷3=5, 0" "93÷3=31, 0"
Replace-Text "25÷7=3, 4" "41÷8=5, 1"
Replace-Text "13÷6=2, 1" "86÷3=28, 2"
Replace-Text "29÷4=7, 1" "28÷9=3, 1"
Replace-Text "48÷6=8, 0" "60÷4=15, 0"
Replace-Text "71÷4=17, 3" "22÷8=2, 6"
Replace-Text "26÷7=3, 5" "96÷9=10, 6"
Replace-Text "88÷7=12, 4" "17÷7=2, 3"
Replace-Text "84÷4=21, 0" "11÷4=2, 3"
Replace-Text "63÷2=31, 1" "58÷7=8, 2"
Replace-Text "29÷7=4, 1" "44÷5=8, 4"
Replace-Text "11÷7=1, 4" "76÷6=12, 4"
Replace-Text "57÷3=19, 0" "20÷8=2, 4"
Replace-Text "63÷8=7, 7" "98÷3=32, 2"
Replace-Text "34÷3=11, 1" "89÷7=12, 5"
Replace-Text "66÷4=16, 2" "68÷4=17, 0"
Replace-Text "87÷9=9, 6" "69÷4=17, 1"
Replace-Text "17÷2=8, 1" "67÷5=13, 2"
Replace-Text "98÷2=49, 0" "91÷6=15, 1"
Replace-Text "91÷5=18, 1" "34÷9=3, 7"
Replace-Text "36÷2=18, 0" "41÷6=6, 5"
Replace-Text "80÷9=8, 8" "80÷7=11, 3"
